$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B4 currently holds the text "-20" (inline string). The parser now reads
# the <input> tag's value as a real number, so store it as a true numeric
# value instead of text.
$ws.Range("B4").Value = -20

# A new row was parsed from the HTML <input> tag: "pawan pritam" / "-12".
# The age value here stays textual (it's still being captured from the
# parsed HTML as a string), so force text storage before assigning it,
# then reset the cell style so only the value semantics differ.
$ws.Range("A5").Value = "pawan pritam"
$ws.Range("B5").NumberFormat = "@"
$ws.Range("B5").Value = "-12"
$ws.Range("B5").Style = "Normal"
